$wb = $excel.ActiveWorkbook

# The localization status text changed from "Ready for handoff" to
# "In Translation". That shared string is used on the Overview sheet
# (columns for zh-cn / de-de) as well as on each language sheet's
# "Status" column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# With the shorter text in place, the status columns were re-sized
# (narrower) to fit the new content.
$wsOverview.Columns.Item(5).ColumnWidth = 12.46
$wsOverview.Columns.Item(6).ColumnWidth = 12.46
$wsZhCn.Columns.Item(3).ColumnWidth = 12.46
$wsDeDe.Columns.Item(3).ColumnWidth = 12.46
